{"js": "// Update the title date paragraph text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.load(\"text\");\nawait context.sync();\nconst oldTitle = \"2025-06-25 Wednesday\";\nconst newTitle = \"2025-06-26 Thursday\";\nif (titleParagraph.text.indexOf(oldTitle) !== -1) {\n  titleParagraph.insertText(newTitle, Word.InsertLocation.replace);\n} else if (titleParagraph.text.trim().length === 0) {\n  titleParagraph.insertText(newTitle, Word.InsertLocation.replace);\n}\n\n// Update the table of arithmetic problems/answers (20 rows x 5 cols).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\n    \"49+48=97\",\n    \"29+63=92\",\n    \"61-23=38\",\n    \"2+49=51\",\n    \"80-62=18\"\n  ],\n  [\n    \"24+28=52\",\n    \"52+19=71\",\n    \"41-38=3\",\n    \"81-36=45\",\n    \"4+39=43\"\n  ],\n  [\n    \"74+17=91\",\n    \"81-4=77\",\n    \"57+5=62\",\n    \"85+8=93\",\n    \"22-7=15\"\n  ],\n  [\n    \"63-17=46\",\n    \"86+9=95\",\n    \"9+77=86\",\n    \"38+5=43\",\n    \"40-33=7\"\n  ],\n  [\n    \"95-67=28\",\n    \"90-9=81\",\n    \"93-38=55\",\n    \"53-16=37\",\n    \"18+59=77\"\n  ],\n  [\n    \"19+42=61\",\n    \"5+69=74\",\n    \"39+24=63\",\n    \"19+73=92\",\n    \"90-18=72\"\n  ],\n  [\n    \"77-49=28\",\n    \"51-7=44\",\n    \"85+9=94\",\n    \"8+87=95\",\n    \"30-18=12\"\n  ],\n  [\n    \"62-35=27\",\n    \"16+47=63\",\n    \"77+14=91\",\n    \"72-48=24\",\n    \"14+78=92\"\n  ],\n  [\n    \"49+12=61\",\n    \"67+8=75\",\n    \"81-79=2\",\n    \"78-29=49\",\n    \"20-7=13\"\n  ],\n  [\n    \"43+9=52\",\n    \"48+48=96\",\n    \"44-35=9\",\n    \"84-58=26\",\n    \"62-47=15\"\n  ],\n  [\n    \"92-43=49\",\n    \"82-33=49\",\n    \"73-69=4\",\n    \"82-9=73\",\n    \"39+13=52\"\n  ],\n  [\n    \"92-13=79\",\n    \"91-17=74\",\n    \"39+49=88\",\n    \"18+58=76\",\n    \"14+8=22\"\n  ],\n  [\n    \"84-15=69\",\n    \"28+17=45\",\n    \"24+57=81\",\n    \"94-67=27\",\n    \"6+18=24\"\n  ],\n  [\n    \"9+17=26\",\n    \"61-46=15\",\n    \"22+39=61\",\n    \"4+9=13\",\n    \"73-28=45\"\n  ],\n  [\n    \"17+18=35\",\n    \"55+6=61\",\n    \"14+79=93\",\n    \"6+37=43\",\n    \"81-23=58\"\n  ],\n  [\n    \"38+35=73\",\n    \"70-69=1\",\n    \"72-43=29\",\n    \"59+27=86\",\n    \"59+32=91\"\n  ],\n  [\n    \"51-45=6\",\n    \"54+18=72\",\n    \"39+42=81\",\n    \"72-3=69\",\n    \"83-79=4\"\n  ],\n  [\n    \"64-18=46\",\n    \"80-33=47\",\n    \"17+14=31\",\n    \"15+46=61\",\n    \"46+38=84\"\n  ],\n  [\n    \"70-22=48\",\n    \"61-39=22\",\n    \"57+35=92\",\n    \"49+4=53\",\n    \"53+19=72\"\n  ],\n  [\n    \"14+38=52\",\n    \"9+4=13\",\n    \"75-59=16\",\n    \"66+9=75\",\n    \"42-19=23\"\n  ]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title date paragraph.\n$titleParagraph = $d.Paragraphs.Item(1)\nif ($titleParagraph.Range.Text -like '*2025-06-25 Wednesday*') {\n  $titleParagraph.Range.Text = '2025-06-26 Thursday'\n}\n\n# Update the table of arithmetic problems/answers (20 rows x 5 cols).\n$t = $d.Tables.Item(1)\n$newValues = @(\n  @('49+48=97', '29+63=92', '61-23=38', '2+49=51', '80-62=18'),\n  @('24+28=52', '52+19=71', '41-38=3', '81-36=45', '4+39=43'),\n  @('74+17=91', '81-4=77', '57+5=62', '85+8=93', '22-7=15'),\n  @('63-17=46', '86+9=95', '9+77=86', '38+5=43', '40-33=7'),\n  @('95-67=28', '90-9=81', '93-38=55', '53-16=37', '18+59=77'),\n  @('19+42=61', '5+69=74', '39+24=63', '19+73=92', '90-18=72'),\n  @('77-49=28', '51-7=44', '85+9=94', '8+87=95', '30-18=12'),\n  @('62-35=27', '16+47=63', '77+14=91', '72-48=24', '14+78=92'),\n  @('49+12=61', '67+8=75', '81-79=2', '78-29=49', '20-7=13'),\n  @('43+9=52', '48+48=96', '44-35=9', '84-58=26', '62-47=15'),\n  @('92-43=49', '82-33=49', '73-69=4', '82-9=73', '39+13=52'),\n  @('92-13=79', '91-17=74', '39+49=88', '18+58=76', '14+8=22'),\n  @('84-15=69', '28+17=45', '24+57=81', '94-67=27', '6+18=24'),\n  @('9+17=26', '61-46=15', '22+39=61', '4+9=13', '73-28=45'),\n  @('17+18=35', '55+6=61', '14+79=93', '6+37=43', '81-23=58'),\n  @('38+35=73', '70-69=1', '72-43=29', '59+27=86', '59+32=91'),\n  @('51-45=6', '54+18=72', '39+42=81', '72-3=69', '83-79=4'),\n  @('64-18=46', '80-33=47', '17+14=31', '15+46=61', '46+38=84'),\n  @('70-22=48', '61-39=22', '57+35=92', '49+4=53', '53+19=72'),\n  @('14+38=52', '9+4=13', '75-59=16', '66+9=75', '42-19=23')\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Cell($r, $c).Range.Text = $newValues[$r-1][$c-1]\n  }\n}\n"}
